$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update entrepreneur poverty-measure values with the recomputed basket figures ---
$ws.Range("C2").Value = 23.56495
$ws.Range("D2").Value = 6.2510240000000001
$ws.Range("E2").Value = 2.4138310000000001
$ws.Range("C4").Value = 5.2306280000000003
$ws.Range("D4").Value = 1.4537690000000001
$ws.Range("E4").Value = 0.7292073
$ws.Range("C5").Value = 26.36542
$ws.Range("D5").Value = 6.9837790000000002
$ws.Range("E5").Value = 2.6711490000000002
$ws.Range("C7").Value = 14.15935
$ws.Range("D7").Value = 3.9708670000000001
$ws.Range("E7").Value = 1.5857760000000001
$ws.Range("C8").Value = 30.28041
$ws.Range("D8").Value = 7.7507450000000002
$ws.Range("E8").Value = 2.8423630000000002
$ws.Range("C9").Value = 21.565909999999999
$ws.Range("D9").Value = 6.6111709999999997
$ws.Range("E9").Value = 2.7968660000000001
$ws.Range("C10").Value = 25.67445
$ws.Range("D10").Value = 6.2199439999999999
$ws.Range("E10").Value = 2.3491810000000002
$ws.Range("C12").Value = 14.15935
$ws.Range("D12").Value = 3.9708670000000001
$ws.Range("E12").Value = 1.5857760000000001
$ws.Range("C13").Value = 30.28041
$ws.Range("D13").Value = 7.7507450000000002
$ws.Range("E13").Value = 2.8423630000000002
$ws.Range("C14").Value = 21.565909999999999
$ws.Range("D14").Value = 6.6111709999999997
$ws.Range("E14").Value = 2.7968660000000001
$ws.Range("C15").Value = 25.67445
$ws.Range("D15").Value = 6.2199439999999999
$ws.Range("E15").Value = 2.3491810000000002
$ws.Range("C17").Value = 2.0357249999999998
$ws.Range("D17").Value = 0.48102640000000002
$ws.Range("E17").Value = 0.1326764
$ws.Range("C18").Value = 12.93952
$ws.Range("D18").Value = 4.1126880000000003
$ws.Range("E18").Value = 1.824268
$ws.Range("C19").Value = 20.536290000000001
$ws.Range("D19").Value = 5.3008350000000002
$ws.Range("E19").Value = 1.955479
$ws.Range("C20").Value = 27.434979999999999
$ws.Range("D20").Value = 6.9974239999999996
$ws.Range("E20").Value = 2.5224959999999998
$ws.Range("C21").Value = 32.602539999999998
$ws.Range("D21").Value = 8.3655209999999993
$ws.Range("E21").Value = 3.1034009999999999
$ws.Range("C22").Value = 13.99034
$ws.Range("D22").Value = 4.0681710000000004
$ws.Range("E22").Value = 1.768157
$ws.Range("C23").Value = 60.384689999999999
$ws.Range("D23").Value = 21.032139999999998
$ws.Range("E23").Value = 9.2366130000000002
$ws.Range("C24").Value = 14.19328
$ws.Range("D24").Value = 3.4699939999999998
$ws.Range("E24").Value = 1.23567
$ws.Range("C25").Value = 29.47803
$ws.Range("D25").Value = 7.5661139999999998
$ws.Range("E25").Value = 3.0220699999999998
$ws.Range("C26").Value = 21.642869999999998
$ws.Range("D26").Value = 4.7930780000000004
$ws.Range("E26").Value = 1.6359570000000001
$ws.Range("C28").Value = 19.04241
$ws.Range("D28").Value = 5.0572359999999996
$ws.Range("E28").Value = 1.855399
$ws.Range("C29").Value = 4.1480319999999997
$ws.Range("D29").Value = 1.743579
$ws.Range("E29").Value = 1.0329900000000001
$ws.Range("C30").Value = 30.28041
$ws.Range("D30").Value = 7.7507450000000002
$ws.Range("E30").Value = 2.8423630000000002
$ws.Range("C31").Value = 21.565909999999999
$ws.Range("D31").Value = 6.6111709999999997
$ws.Range("E31").Value = 2.7968660000000001
$ws.Range("C32").Value = 25.67445
$ws.Range("D32").Value = 6.2199439999999999
$ws.Range("E32").Value = 2.3491810000000002

# --- Region labels for the 6-region TPI breakdown switch from numeric codes to names ---
$ws.Range("H7").Value = "Central"
$ws.Range("H8").Value = "Eastern"
$ws.Range("H9").Value = "Northern"
$ws.Range("H10").Value = "Western"

# --- Capitalize the region-name labels used alongside the TPI table ---
$ws.Range("I12").Value = "Central"
$ws.Range("I13").Value = "Eastern"
$ws.Range("I14").Value = "Northern"
$ws.Range("I15").Value = "Western"

# --- Spatial domain label: "Central rural" -> "central rural" ---
$ws.Range("K28").Value = "central rural"

# --- Remove the stray AVERAGE() helper row that is no longer part of the write-up ---
$ws.Rows.Item(35).Delete()
